# requerimento-para-desarquivamento.docx
#
# The paragraph that ends with "{{ cidade_assinatura }}, {{ data }}."
# needs its final run split so the placeholder "{{ data }}" becomes
# "{{ data_assinatura }}" - i.e. insert the literal text "_assinatura"
# right after "{{ data" and before " }}.".
#
# The original run (" }}, {{ data }}.") must become three runs:
#   1) " }}, {{ data"      (keeps the original run's formatting/rsid)
#   2) "_assinatura"        (new run, same rPr, no rsid attributes)
#   3) " }}."               (new run, same rPr, no rsid attributes)

$d = $word.ActiveDocument

# Locate the unique run of text that needs to be split.
$found = $d.Content
$found.Find.Execute(" }}, {{ data }}.", $true, $false, $false, $false, $false, `
                     $true, 1, $false, "", 0)

if (-not $found.Find.Found) {
    throw "Could not find the target text ' }}, {{ data }}.' in the document."
}

$rangeStart = $found.Start
$rangeEnd = $found.End
$target = $d.Range($rangeStart, $rangeEnd)

# Re-express that range as three runs, splitting out "_assinatura" as its
# own run (matching how Word splits a run when new text is typed into the
# middle of it). The first run keeps the original rsid attribute; the two
# brand-new runs carry no rsid attributes.
$xmlFragment = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData>' + `
      '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:body>' + `
          '<w:p>' + `
            '<w:r w:rsidRPr="001D437B">' + `
              '<w:rPr>' + `
                '<w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' + `
                '<w:sz w:val="24"/>' + `
                '<w:szCs w:val="24"/>' + `
              '</w:rPr>' + `
              '<w:t xml:space="preserve"> }}, {{ data</w:t>' + `
            '</w:r>' + `
            '<w:r>' + `
              '<w:rPr>' + `
                '<w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' + `
                '<w:sz w:val="24"/>' + `
                '<w:szCs w:val="24"/>' + `
              '</w:rPr>' + `
              '<w:t>_assinatura</w:t>' + `
            '</w:r>' + `
            '<w:r>' + `
              '<w:rPr>' + `
                '<w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' + `
                '<w:sz w:val="24"/>' + `
                '<w:szCs w:val="24"/>' + `
              '</w:rPr>' + `
              '<w:t xml:space="preserve"> }}.</w:t>' + `
            '</w:r>' + `
          '</w:p>' + `
        '</w:body>' + `
      '</w:document>' + `
    '</pkg:xmlData>' + `
  '</pkg:part>' + `
'</pkg:package>'

$target.InsertXML($xmlFragment)
